$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.788.54'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.619.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.991'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.78%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.77'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.521'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.990'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.29'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +8.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.259'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.34%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.850.34'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.620.94'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.782.19'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.94'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +17.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.37'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.41%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.09'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.94'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.64'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.39%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.79%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.992'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.96%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.34'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.421.33'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.63'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.44%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.88'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.56%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.59%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.97%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.827'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.53%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '53.68'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.43'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.74%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +18.13%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.760.21'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.20'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0532'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.97%  '
Write-Host "Applied cryptos price/volume update"
